$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C, shifting existing C.. columns to the right
$ws.Columns.Item(3).Insert() | Out-Null

# Set width of new column C to match column B's display width, not best-fit
$ws.Columns.Item(3).ColumnWidth = 21.3

# Populate the new "Pan *" column
$ws.Range("C1").Value = "Pan *"
$ws.Range("C2").Value = "BUHNXDFEA7"
$ws.Range("C3").Value = "BUHNXDFEA7"

# Update the selection to match the target state
$ws.Range("C4").Select() | Out-Null
